$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 0.971099
$ws.Range("H2").Value = 2.913297
$ws.Range("I2").Value = 0.007882574716876797
$ws.Range("J2").Value = 0.007882574716876797
$ws.Range("M2").Value = 0.967553
$ws.Range("N2").Value = 2.902659
$ws.Range("O2").Value = 0.01359591950841534
$ws.Range("P2").Value = 0.01359591950841534
$ws.Range("Q2").Value = 0.939589750747
$ws.Range("R2").Value = 8.456307756723
$ws.Range("S2").Value = 0.0001071708513697268
$ws.Range("T2").Value = 0.0001071708513697268
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 0.971099
$ws.Range("H3").Value = 2.913297
$ws.Range("I3").Value = 0.007882574716876797
$ws.Range("J3").Value = 0.007882574716876797
$ws.Range("M3").Value = 5.176377666666667
$ws.Range("O3").Value = 0.07273773540173906
$ws.Range("P3").Value = 0.07273773540173906
$ws.Range("Q3").Value = 5.026775175722334
$ws.Range("R3").Value = 45.240976581501
$ws.Range("S3").Value = 0.0005733606340406226
$ws.Range("T3").Value = 0.0005733606340406226
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 0.971099
$ws.Range("H4").Value = 2.913297
$ws.Range("I4").Value = 0.007882574716876797
$ws.Range("J4").Value = 0.007882574716876797
$ws.Range("M4").Value = 4.670153
$ws.Range("N4").Value = 14.010459
$ws.Range("O4").Value = 0.06562433714740633
$ws.Range("P4").Value = 0.06562433714740633
$ws.Range("Q4").Value = 4.535180908147
$ws.Range("R4").Value = 40.816628173323
$ws.Range("S4").Value = 0.000517288740809944
$ws.Range("T4").Value = 0.000517288740809944
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 0.971099
$ws.Range("H5").Value = 2.913297
$ws.Range("I5").Value = 0.007882574716876797
$ws.Range("J5").Value = 0.007882574716876797
$ws.Range("M5").Value = 60.35087133333334
$ws.Range("N5").Value = 181.052614
$ws.Range("O5").Value = 0.8480420079424392
$ws.Range("P5").Value = 0.8480420079424392
$ws.Range("Q5").Value = 58.60667080092868
$ws.Range("R5").Value = 527.460037208358
$ws.Range("S5").Value = 0.006684754490656503
$ws.Range("T5").Value = 0.006684754490656503
$ws.Range("I6").Value = 0.6966643430097871
$ws.Range("J6").Value = 0.696664343009787
$ws.Range("M6").Value = 0.967553
$ws.Range("N6").Value = 2.902659
$ws.Range("O6").Value = 0.01359591950841534
$ws.Range("P6").Value = 0.01359591950841534
$ws.Range("Q6").Value = 83.041226999272
$ws.Range("R6").Value = 747.371042993448
$ws.Range("S6").Value = 0.009471792331944121
$ws.Range("T6").Value = 0.009471792331944121
$ws.Range("I7").Value = 0.6966643430097871
$ws.Range("J7").Value = 0.696664343009787
$ws.Range("M7").Value = 5.176377666666667
$ws.Range("O7").Value = 0.07273773540173906
$ws.Range("P7").Value = 0.07273773540173906
$ws.Range("Q7").Value = 444.2679138523974
$ws.Range("S7").Value = 0.05067378664567227
$ws.Range("T7").Value = 0.05067378664567226
$ws.Range("I8").Value = 0.6966643430097871
$ws.Range("J8").Value = 0.696664343009787
$ws.Range("M8").Value = 4.670153
$ws.Range("N8").Value = 14.010459
$ws.Range("O8").Value = 0.06562433714740633
$ws.Range("P8").Value = 0.06562433714740633
$ws.Range("Q8").Value = 400.820663461672
$ws.Range("R8").Value = 3607.385971155048
$ws.Range("S8").Value = 0.0457181357242506
$ws.Range("T8").Value = 0.04571813572425059
$ws.Range("I9").Value = 0.6966643430097871
$ws.Range("J9").Value = 0.696664343009787
$ws.Range("M9").Value = 60.35087133333334
$ws.Range("N9").Value = 181.052614
$ws.Range("O9").Value = 0.8480420079424392
$ws.Range("P9").Value = 0.8480420079424392
$ws.Range("Q9").Value = 5179.67533147558
$ws.Range("R9").Value = 46617.07798328021
$ws.Range("S9").Value = 0.59080062830792
$ws.Range("T9").Value = 0.59080062830792
$ws.Range("G10").Value = 36.24916566666667
$ws.Range("H10").Value = 108.747497
$ws.Range("I10").Value = 0.294240604502677
$ws.Range("J10").Value = 0.294240604502677
$ws.Range("M10").Value = 0.967553
$ws.Range("N10").Value = 2.902659
$ws.Range("O10").Value = 0.01359591950841534
$ws.Range("P10").Value = 0.01359591950841534
$ws.Range("Q10").Value = 35.07298898828034
$ws.Range("R10").Value = 315.656900894523
$ws.Range("S10").Value = 0.00400047157492587
$ws.Range("T10").Value = 0.004000471574925869
$ws.Range("G11").Value = 36.24916566666667
$ws.Range("H11").Value = 108.747497
$ws.Range("I11").Value = 0.294240604502677
$ws.Range("J11").Value = 0.294240604502677
$ws.Range("M11").Value = 5.176377666666667
$ws.Range("O11").Value = 0.07273773540173906
$ws.Range("P11").Value = 0.07273773540173906
$ws.Range("Q11").Value = 187.6393715922335
$ws.Range("R11").Value = 1688.754344330101
$ws.Range("S11").Value = 0.02140239523476347
$ws.Range("T11").Value = 0.02140239523476347
$ws.Range("G12").Value = 36.24916566666667
$ws.Range("H12").Value = 108.747497
$ws.Range("I12").Value = 0.294240604502677
$ws.Range("J12").Value = 0.294240604502677
$ws.Range("M12").Value = 4.670153
$ws.Range("N12").Value = 14.010459
$ws.Range("O12").Value = 0.06562433714740633
$ws.Range("P12").Value = 0.06562433714740633
$ws.Range("Q12").Value = 169.2891497856803
$ws.Range("R12").Value = 1523.602348071123
$ws.Range("S12").Value = 0.01930934463234032
$ws.Range("T12").Value = 0.01930934463234032
$ws.Range("G13").Value = 36.24916566666667
$ws.Range("H13").Value = 108.747497
$ws.Range("I13").Value = 0.294240604502677
$ws.Range("J13").Value = 0.294240604502677
$ws.Range("M13").Value = 60.35087133333334
$ws.Range("N13").Value = 181.052614
$ws.Range("O13").Value = 0.8480420079424392
$ws.Range("P13").Value = 0.8480420079424392
$ws.Range("Q13").Value = 2187.668733089685
$ws.Range("R13").Value = 19689.01859780716
$ws.Range("S13").Value = 0.2495283930606474
$ws.Range("T13").Value = 0.2495283930606473
$ws.Range("E14").Value = 2.0
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.149372
$ws.Range("H14").Value = 0.448116
$ws.Range("I14").Value = 0.001212477770659141
$ws.Range("J14").Value = 0.001212477770659141
$ws.Range("M14").Value = 0.967553
$ws.Range("N14").Value = 2.902659
$ws.Range("O14").Value = 0.01359591950841534
$ws.Range("P14").Value = 0.01359591950841534
$ws.Range("Q14").Value = 0.144525326716
$ws.Range("R14").Value = 1.300727940444
$ws.Range("S14").Value = 0.00001648475017562455
$ws.Range("T14").Value = 0.00001648475017562455
$ws.Range("E15").Value = 2.0
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.149372
$ws.Range("H15").Value = 0.448116
$ws.Range("I15").Value = 0.001212477770659141
$ws.Range("J15").Value = 0.001212477770659141
$ws.Range("M15").Value = 5.176377666666667
$ws.Range("O15").Value = 0.07273773540173906
$ws.Range("P15").Value = 0.07273773540173906
$ws.Range("Q15").Value = 0.7732058848253333
$ws.Range("R15").Value = 6.958852963428
$ws.Range("S15").Value = 0.00008819288726269502
$ws.Range("T15").Value = 0.00008819288726269502
$ws.Range("E16").Value = 2.0
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.149372
$ws.Range("H16").Value = 0.448116
$ws.Range("I16").Value = 0.001212477770659141
$ws.Range("J16").Value = 0.001212477770659141
$ws.Range("M16").Value = 4.670153
$ws.Range("N16").Value = 14.010459
$ws.Range("O16").Value = 0.06562433714740633
$ws.Range("P16").Value = 0.06562433714740633
$ws.Range("Q16").Value = 0.6975900939159999
$ws.Range("R16").Value = 6.278310845243999
$ws.Range("S16").Value = 0.00007956805000547106
$ws.Range("T16").Value = 0.00007956805000547106
$ws.Range("E17").Value = 2.0
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.149372
$ws.Range("H17").Value = 0.448116
$ws.Range("I17").Value = 0.001212477770659141
$ws.Range("J17").Value = 0.001212477770659141
$ws.Range("M17").Value = 60.35087133333334
$ws.Range("N17").Value = 181.052614
$ws.Range("O17").Value = 0.8480420079424392
$ws.Range("P17").Value = 0.8480420079424392
$ws.Range("Q17").Value = 9.014730352802665
$ws.Range("R17").Value = 81.132573175224
$ws.Range("S17").Value = 0.00102823208321535
$ws.Range("T17").Value = 0.00102823208321535
